# Web App and Script Access Update
# Slide 4 ("How it works"), content placeholder: the second bullet's
# label "Authentication:" loses its trailing colon, becoming
# "Authentication" (formatting/run properties are left untouched).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(2)
$para.Text = "Authentication"
